$d = $word.ActiveDocument

# --- Step 1: Add the SubtleEmphasis character style (used later). ---
$subtle = $d.Styles.Add("Subtle Emphasis", 2)
$subtleFont = $subtle.Font
$subtleFont.Italic = $true
$subtleFont.ItalicBi = $true
$subtleTc = $subtleFont.TextColor
$subtleTc.ObjectThemeColor = 13
$subtle.BaseStyle = "DefaultParagraphFont"
$subtle.Priority = 19
$subtle.QuickStyle = $true

# --- Step 2: Convert the Task A.5 field-code hyperlink into a real w:hyperlink. ---
# (paragraph 7 at this point: "Task A.5")
$fields = $d.Fields
$task5Field = $fields.Item($fields.Count)
$task5Field.Unlink()

$p5 = $d.Paragraphs.Item(7)
$p5Range = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$d.Hyperlinks.Add($p5Range, "", "_Task_A.5") | Out-Null

# --- Step 3: Center-align the Task A.1 - Task A.7 paragraphs (indices 3-9). ---
for ($i = 3; $i -le 9; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Format.Alignment = 1
}

# --- Step 4: Insert the new "Task A" paragraph right after "Table of Contents". ---
$toc = $d.Paragraphs.Item(2)
$toc.Range.InsertParagraphAfter()
$taskAHeader = $d.Paragraphs.Item(3)
$taskAHeader.Style = "Normal"
$taskAHeader.Range.Text = "Task A"
$taskAHeader.Format.Alignment = 1

# --- Step 5: Insert the "Ctrl + Left-click..." paragraph and an empty paragraph after Task A.7. ---
# Task A.7 is now at index 10 (shifted by the insert in step 4).
$taskA7 = $d.Paragraphs.Item(10)
$taskA7.Range.InsertParagraphAfter()
$notePara = $d.Paragraphs.Item(11)
$noteXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rStyle w:val="SubtleEmphasis"/></w:rPr><w:t>Ctrl + Left-click to follow hyperlinks.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$notePara.Range.InsertXML($noteXml)
$notePara.Format.Alignment = 1

$notePara.Range.InsertParagraphAfter()
$emptyPara = $d.Paragraphs.Item(12)
$emptyPara.Format.Alignment = 1

# --- Step 6: Change the "Task A" section heading from Subtitle+center to Heading2. ---
$sectionHeading = $d.Paragraphs.Item(14)
$sectionHeading.Style = "Heading2"
$sectionHeading.Format.Alignment = 0

Write-Output "done"
